# Surveying Phylogenetic Forests - apply authored edits
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) "Future Steps?" slide (slide 19): fill in the previously empty
#    content placeholder with three new bullet lines.
# ---------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$tf19 = $s19.Shapes.Item(2).TextFrame
$tr19 = $tf19.TextRange
$tr19.Text = "Exploring compression methods"

$tr19b = $tf19.TextRange
$tr19b.InsertAfter("`r7zip vs. WinZip") | Out-Null

$tr19c = $tf19.TextRange
$tr19c.InsertAfter("`rBiological assumptions in compression (MFC)") | Out-Null

$tr19d = $tf19.TextRange
$tr19d.Paragraphs(2,1).IndentLevel = 2
$tr19d.Paragraphs(3,1).IndentLevel = 2

# ---------------------------------------------------------------
# 2) "Draft" slide (slide 2): fix the "liklihood" typo and merge the
#    "Maximum " / "liklihood" runs into a single "Maximum likelihood" run.
# ---------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tf2 = $s2.Shapes.Item(2).TextFrame
$tr2 = $tf2.TextRange
$paraCount2 = $tr2.Paragraphs().Count
for ($i = 1; $i -le $paraCount2; $i++) {
  $para = $tr2.Paragraphs($i,1)
  if ($para.Text -like "*liklihood*") {
    $sub = $tr2.Characters($para.Start, $para.Length)
    $sub.Text = "Maximum likelihood"
  }
}

# ---------------------------------------------------------------
# 3) "MrBayes" slide (slide 5): add a new red sub-bullet "costly
#    precursor" right after "Multiple Sequence Alignment".
# ---------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tf5 = $s5.Shapes.Item(2).TextFrame
$tr5 = $tf5.TextRange
$paraCount5 = $tr5.Paragraphs().Count
for ($i = 1; $i -le $paraCount5; $i++) {
  $para = $tr5.Paragraphs($i,1)
  if ($para.Text -like "*Multiple Sequence Alignment*") {
    $ins = $para.InsertAfter("`rcostly precursor")
  }
}

$tr5b = $tf5.TextRange
$paraCount5b = $tr5b.Paragraphs().Count
for ($i = 1; $i -le $paraCount5b; $i++) {
  $para = $tr5b.Paragraphs($i,1)
  if ($para.Text -like "*costly precursor*") {
    $para.IndentLevel = 2
    $para.Font.Color.RGB = 255
  }
}

# ---------------------------------------------------------------
# 4) "Normalized Compression Distance" slide (slide 6): add a new red
#    sub-bullet "Muscle timings" right after "Avoids Multiple Sequence
#    Alignment".
# ---------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tf6 = $s6.Shapes.Item(2).TextFrame
$tr6 = $tf6.TextRange
$paraCount6 = $tr6.Paragraphs().Count
for ($i = 1; $i -le $paraCount6; $i++) {
  $para = $tr6.Paragraphs($i,1)
  if ($para.Text -like "*Avoids Multiple Sequence Alignment*") {
    $ins = $para.InsertAfter("`rMuscle timings")
  }
}

$tr6b = $tf6.TextRange
$paraCount6b = $tr6b.Paragraphs().Count
for ($i = 1; $i -le $paraCount6b; $i++) {
  $para = $tr6b.Paragraphs($i,1)
  if ($para.Text -like "*Muscle timings*") {
    $para.IndentLevel = 2
    $para.Font.Color.RGB = 255
  }
}
